$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @("航天发展", "航天发展", "航天发展"),
    @("利欧股份", "白银有色", "利欧股份"),
    @("浙文互联", "利欧股份", "白银有色"),
    @("巨力索具", "湖南黄金", "铜陵有色"),
    @("白银有色", "浙文互联", "巨力索具"),
    @("中超控股", "巨力索具", "浙文互联"),
    @("蓝色光标", "铜陵有色", "湖南黄金"),
    @("特变电工", "中超控股", "湖南白银"),
    @("湖南黄金", "特变电工", "中超控股"),
    @("铜陵有色", "红 宝 丽", "洲际油气"),
    @("信维通信", "湖南白银", "红宝丽"),
    @("湖南白银", "蓝色光标", "华胜天成"),
    @("红 宝 丽", "华胜天成", "蓝色光标"),
    @("网宿科技", "紫金矿业", "西部材料"),
    @("通宇通讯", "寒武纪-U", "网宿科技"),
    @("华胜天成", "中国黄金", "紫金矿业"),
    @("紫金矿业", "通宇通讯", "晓程科技"),
    @("西部材料", "岩山科技", "航天电子"),
    @("海兰信", "航天电子", "通宇通讯"),
    @("三变科技", "信维通信", "锋龙股份")
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
    $ws.Cells.Item($row, 3).Value = $values[$i][2]
}
